$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.340.65'
$ws.Range("E2").Value = '  +4.18%  '

$ws.Range("D3").Value = '1.716.01'
$ws.Range("E3").Value = '  +3.29%  '

$ws.Range("D4").Value = '''0.9987'
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '''239.93'
$ws.Range("E5").Value = '  +1.73%  '

$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Value = '''0.4716'
$ws.Range("E7").Value = '  -1.89%  '

$ws.Range("D8").Value = '''0.2634'
$ws.Range("E8").Value = '  +1.52%  '

$ws.Range("D9").Value = '''0.06226'
$ws.Range("E9").Value = '  +1.33%  '

$ws.Range("D10").Value = '1.712.77'
$ws.Range("E10").Value = '  +3.51%  '

$ws.Range("D11").Value = '''0.07079'
$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("D12").Value = '''15.23'
$ws.Range("E12").Value = '  +3.82%  '

$ws.Range("D13").Value = '''0.5904'
$ws.Range("E13").Value = '  +0.77%  '

$ws.Range("D14").Value = '''4.419'
$ws.Range("E14").Value = '  +1.31%  '

$ws.Range("D15").Value = '''76.15'
$ws.Range("E15").Value = '  +2.36%  '

$ws.Range("D16").Value = '''0.9996'
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '''0.9998'
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").Value = '26.330.14'
$ws.Range("E18").Value = '  +4.21%  '

$ws.Range("D19").Value = '''0.000006795'
$ws.Range("E19").Value = '  +1.60%  '

$ws.Range("D20").Value = '''11.58'
$ws.Range("E20").Value = '  +1.80%  '

$ws.Range("D21").Value = '1.926.14'
$ws.Range("E21").Value = '  +3.60%  '

$ws.Range("D22").Value = '''4.565'
$ws.Range("E22").Value = '  +4.67%  '

$ws.Range("D23").Value = '''8.846'
$ws.Range("E23").Value = '  +2.92%  '

$ws.Range("D24").Value = '''5.343'
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").Value = '''135.71'
$ws.Range("E25").Value = '  +0.87%  '

$ws.Range("D26").Value = '''15.18'
$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("D27").Value = '''1.405'
$ws.Range("E27").Value = '  +1.97%  '

$ws.Range("D28").Value = '''1.767'
$ws.Range("E28").Value = '  +5.27%  '

$ws.Range("D29").Value = '''106.58'
$ws.Range("E29").Value = '  +1.78%  '

$ws.Range("D30").Value = '''4.031'
$ws.Range("E30").Value = '  +1.94%  '

$ws.Range("D31").Value = '''3.686'
$ws.Range("E31").Value = '  +1.96%  '

$ws.Range("D32").Value = '''0.07767'
$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("D33").Value = '''0.04422'
$ws.Range("E33").Value = '  +1.96%  '

$ws.Range("D34").Value = '''2.611'
$ws.Range("E34").Value = '  +0.84%  '

$ws.Range("D35").Value = '''0.6233'
$ws.Range("E35").Value = '  +3.69%  '

$ws.Range("D36").Value = '''0.9719'
$ws.Range("E36").Value = '  +3.07%  '

$ws.Range("D37").Value = '''0.9156'
$ws.Range("E37").Value = '  +8.00%  '

$ws.Range("D38").Value = '''112.25'
$ws.Range("E38").Value = '  +13.07%  '

$ws.Range("D39").Value = '''2.404'
$ws.Range("E39").Value = '  -7.85%  '

$ws.Range("D40").Value = '''1.001'
$ws.Range("E40").Value = '  +0.25%  '

$ws.Range("E41").Value = '  +5.41%  '

$ws.Range("D42").Value = '''0.01467'
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.199'
$ws.Range("E43").Value = '  +11.67%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.3813'
$ws.Range("E44").Value = '  +2.06%  '

$ws.Range("E45").Value = '  +3.17%  '

$ws.Range("D46").Value = '''6.253'
$ws.Range("E46").Value = '  +0.79%  '

$ws.Range("D47").Value = '''0.05297'
$ws.Range("E47").Value = '  +0.97%  '

$ws.Range("D48").Value = '''30.65'
$ws.Range("E48").Value = '  +4.13%  '

$ws.Range("D49").Value = '''7.671'
$ws.Range("E49").Value = '  +5.46%  '

$ws.Range("D50").Value = '''1.223'
$ws.Range("E50").Value = '  +0.58%  '

$ws.Range("D51").Value = '''0.3384'
$ws.Range("E51").Value = '  +1.75%  '
